$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(1.101058557113987, 1.077114746702648, 1.062865150140539, 1.057172623236767, 1.056234305724388, 1.062787915088734, 1.09270929310162, 1.154947221653202, 1.202818254818368, 1.225056426734824, 1.233543220176443, 1.23171252684935, 1.225753327833928, 1.222111680058276, 1.201374165388387, 1.188770071696183, 1.181563999359213, 1.179131633284612, 1.190107304949862, 1.227501911817598, 1.252323974423689, 1.239041192170276, 1.189502616561185, 1.137731624592107)
$Dvals = @(0.192460773227225, 0.1921783340519028, 0.1920649250370872, 0.192033856226395, 0.1920296137301065, 0.192064444630283, 0.1923509665291476, 0.193386681095312, 0.1944336422787956, 0.1949714521842765, 0.1951839030125058, 0.1951377577029305, 0.1949887547180253, 0.1948986296583826, 0.1943997287487491, 0.1941093918659931, 0.1939481937256886, 0.1938946119666696, 0.1941396992683195, 0.1950322822868884, 0.1956668685959002, 0.1953235086370597, 0.1941259794702717, 0.1930560591717239)
$Evals = @(0.2103385038075398, 0.211461293987087, 0.2122080888849247, 0.2125268763087274, 0.2125806850795371, 0.2122123295662046, 0.2107137467469826, 0.2082291545515833, 0.2066788228956149, 0.2060329035485591, 0.205796813661614, 0.2058472819674364, 0.20601330997461, 0.2061161138398298, 0.2067222269383482, 0.2071092357745403, 0.2073374194212203, 0.2074156387603363, 0.2070674600440032, 0.2059643128909094, 0.2052929082909465, 0.2056467229934871, 0.2070863291355529, 0.2088528666014575)
$Fvals = @(1.723840628001696, 1.728848985728661, 1.732691006530068, 1.734449707927617, 1.73475340472497, 1.732713943073492, 1.725408470272846, 1.717159150954238, 1.714793836201864, 1.714518441189966, 1.714529104272842, 1.714521698192073, 1.714517015702839, 1.714529112100536, 1.714827928967111, 1.715216170780494, 1.715514843838832, 1.715628919058403, 1.715167042885852, 1.714515272777277, 1.714759233571968, 1.714567786570214, 1.715189018516753, 1.718741195572818)
$Gvals = @(0.002484740335981277, 0.00248726880107774, 0.00248890586971608, 0.002489594319623674, 0.002489709926327893, 0.002488915067879085, 0.002485594633567351, 0.002479751531908843, 0.002475862001265224, 0.002474179297939689, 0.0024735544999326, 0.002473688510600794, 0.002474127647089253, 0.002474398244809264, 0.002475973708950071, 0.00247696236111863, 0.002477539168300346, 0.002477735868511743, 0.002476856273002759, 0.002473998325854741, 0.002472202773565448, 0.002473154498589767, 0.002476904209101451, 0.002481261118532824)
$Jvals = @(0.2777475623083427, 0.2796013600975842, 0.2808065774977031, 0.2813145828211603, 0.2813999563943304, 0.2808133602943972, 0.2783728670062553, 0.2741172431735688, 0.2713120466771768, 0.2701052979079916, 0.2696582779790919, 0.2697541095652607, 0.2700683220330076, 0.2702620813086565, 0.2713923040865636, 0.2721034049414497, 0.2725189401644812, 0.2726607552761546, 0.2720270314983431, 0.2699757603642707, 0.2686931198179447, 0.2693723911511192, 0.2720615390175158, 0.275211929876777)
$Kvals = @(0.4231933473978984, 0.3767427502663452, 0.3482380443062141, 0.336626867931983, 0.3346991456776607, 0.3480814319217131, 0.4071742128405731, 0.5231596533615175, 0.6084147525859294, 0.6472038826792357, 0.6618926360590933, 0.6587291551472845, 0.6484123364260483, 0.6420929850297625, 0.6058798528996192, 0.5836653244293473, 0.570888738700063, 0.5665629357445425, 0.5860300396935259, 0.6514426379026759, 0.6941940952922323, 0.6713770307795528, 0.5849609679323748, 0.4917737115941065)
$Lvals = @(0.1702510204662246, 0.1611192878145999, 0.1555667541306462, 0.1533178900329801, 0.1529453084516632, 0.1555363689159606, 0.1670912053546374, 0.1901758310244617, 0.2073894388844479, 0.2152741686547301, 0.2182675593238343, 0.2176225432162653, 0.215520285103139, 0.2142335777492264, 0.2068752283608006, 0.2023748696174863, 0.1997914918287904, 0.1989176874090219, 0.2028534129417778, 0.2161375640100687, 0.2248638270795027, 0.2202024605539776, 0.2026370512245848, 0.1838858965009535)
$Nvals = @(2.22843271254706, 2.249179115810777, 2.262571177744029, 2.268192993887338, 2.269136426100768, 2.262646329612979, 2.235450442124066, 2.187302296827693, 2.155081141790747, 2.141106259372561, 2.135912441962085, 2.137026657449394, 2.140676993469141, 2.142925713960435, 2.156008170330137, 2.164208682819683, 2.168989657926145, 2.170619450291734, 2.163329073904591, 2.139602137061474, 2.124667407556341, 2.132585997620291, 2.163726538288002, 2.199773448782866)
$Ovals = @(4.224683834273407, 4.239664043588675, 4.250827949614802, 4.255872057541922, 4.256739518643116, 4.250893972619139, 4.229441177480766, 4.202958732303955, 4.192990986065155, 4.190514408256661, 4.189872235437917, 4.189997393647729, 4.190455652441329, 4.190774844403848, 4.193194222586158, 4.195205302715351, 4.196555727945622, 4.19704623556828, 4.194971173896505, 4.19031302881308, 4.188991837701053, 4.18953939600209, 4.195076418503163, 4.208455648071009)

$rows = 2..25
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $Bvals[$i]
    $ws.Cells.Item($r, 4).Value = $Dvals[$i]
    $ws.Cells.Item($r, 5).Value = $Evals[$i]
    $ws.Cells.Item($r, 6).Value = $Fvals[$i]
    $ws.Cells.Item($r, 7).Value = $Gvals[$i]
    $ws.Cells.Item($r, 10).Value = $Jvals[$i]
    $ws.Cells.Item($r, 11).Value = $Kvals[$i]
    $ws.Cells.Item($r, 12).Value = $Lvals[$i]
    $ws.Cells.Item($r, 14).Value = $Nvals[$i]
    $ws.Cells.Item($r, 15).Value = $Ovals[$i]
}
